$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 74 (Celine, 2020-12)
$ws.Range("D74").Value = 1429
$ws.Range("E74").Value = 51919
$ws.Range("F74").Value = 36.33240027991602
$ws.Range("G74").Value = 30

# Row 75 (Sam, 2020-12)
$ws.Range("D75").Value = 1823
$ws.Range("E75").Value = 69102
$ws.Range("F75").Value = 37.90565002742732
$ws.Range("G75").Value = 28
